# CompStat weekly update: new crime data collected.
# Updates the report week/volume header text and the precinct crime-stat
# table (rows 14-29) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: volume/number and the reporting week dates.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# ---------------------------------------------------------------------
# Helper: write a text value into a cell that otherwise holds numbers,
# matching the workbook's existing convention of using the literal
# strings "0" / "***.*" in place of blank/undefined numeric figures.
# Forcing the NumberFormat to Text ("@") before assignment keeps Excel
# from re-parsing the literal "0" back into a number; re-applying the
# formatting (not the value) from a cell that already uses this
# convention keeps the visual style consistent with its neighbors.
# ---------------------------------------------------------------------
$fmtDonor = $ws.Range("D14")

function Set-TextCell($addr) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
}

function Apply-DonorFormat($addr) {
    $fmtDonor.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Row 14 - Murder
$ws.Range("F14").Value = 2
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = -77.777777777777

# Row 15 - Rape
Set-TextCell "C15"
$ws.Range("C15").Value = "0"
Apply-DonorFormat "C15"

Set-TextCell "G15"
$ws.Range("G15").Value = "0"
Apply-DonorFormat "G15"

Set-TextCell "H15"
$ws.Range("H15").Value = "***.*"
Apply-DonorFormat "H15"

# Row 16 - Robbery
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -36.363636363636
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -7.407407407407
$ws.Range("I16").Value = 165
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = 10
$ws.Range("L16").Value = 70.103092783505
$ws.Range("M16").Value = 20.437956204379
$ws.Range("N16").Value = -63.087248322147

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 55
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = 37.5
$ws.Range("I17").Value = 243
$ws.Range("J17").Value = 226
$ws.Range("K17").Value = 7.522123893805
$ws.Range("L17").Value = 49.079754601227
$ws.Range("M17").Value = 16.826923076923
$ws.Range("N17").Value = -6.538461538461

# Row 18 - Burglary
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 97
$ws.Range("J18").Value = 81
$ws.Range("K18").Value = 19.753086419753
$ws.Range("L18").Value = 97.959183673469
$ws.Range("M18").Value = -15.652173913043
$ws.Range("N18").Value = -81.593927893738

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -21.153846153846
$ws.Range("I19").Value = 218
$ws.Range("J19").Value = 237
$ws.Range("K19").Value = -8.016877637130
$ws.Range("L19").Value = 37.106918238993
$ws.Range("M19").Value = 52.447552447552
$ws.Range("N19").Value = -7.627118644067

# Row 20 - G.L.A.
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 175
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 36.842105263157
$ws.Range("I20").Value = 181
$ws.Range("J20").Value = 126
$ws.Range("K20").Value = 43.650793650793
$ws.Range("L20").Value = 132.051282051282
$ws.Range("M20").Value = 341.463414634146
$ws.Range("N20").Value = -17.351598173516

# Row 21 - TOTAL
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 160
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = 11.111111111111
$ws.Range("I21").Value = 926
$ws.Range("J21").Value = 839
$ws.Range("K21").Value = 10.369487485101
$ws.Range("L21").Value = 61.888111888111
$ws.Range("M21").Value = 40.943683409436
$ws.Range("N21").Value = -46.535796766743

# Row 22 - Transit
Set-TextCell "C22"
$ws.Range("C22").Value = "0"
Apply-DonorFormat "C22"

Set-TextCell "D22"
$ws.Range("D22").Value = "0"
Apply-DonorFormat "D22"

Set-TextCell "E22"
$ws.Range("E22").Value = "***.*"
Apply-DonorFormat "E22"

$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("M22").Value = -7.692307692307

# Row 23 - Housing
Set-TextCell "C23"
$ws.Range("C23").Value = "0"
Apply-DonorFormat "C23"

Set-TextCell "D23"
$ws.Range("D23").Value = "0"
Apply-DonorFormat "D23"

Set-TextCell "E23"
$ws.Range("E23").Value = "***.*"
Apply-DonorFormat "E23"

$ws.Range("M23").Value = 133.333333333333

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 113.333333333333
$ws.Range("F24").Value = 80
$ws.Range("G24").Value = 62
$ws.Range("H24").Value = 29.032258064516
$ws.Range("I24").Value = 436
$ws.Range("J24").Value = 450
$ws.Range("K24").Value = -3.111111111111
$ws.Range("L24").Value = 68.992248062015
$ws.Range("M24").Value = 43.894389438943

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 77.777777777777
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = -12.727272727272
$ws.Range("I25").Value = 283
$ws.Range("J25").Value = 309
$ws.Range("K25").Value = -8.414239482200
$ws.Range("L25").Value = 41.5
$ws.Range("M25").Value = -18.678160919540

# Row 26 - UCR Rape*
Set-TextCell "C26"
$ws.Range("C26").Value = "0"
Apply-DonorFormat "C26"

$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50

# Row 27 - Other Sex Crimes
Set-TextCell "C27"
$ws.Range("C27").Value = "0"
Apply-DonorFormat "C27"

$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 9
$ws.Range("H27").Value = 80
$ws.Range("J27").Value = 62
$ws.Range("K27").Value = -11.290322580645
$ws.Range("L27").Value = -17.910447761194

# Row 28 - Shooting Vic.
Set-TextCell "C28"
$ws.Range("C28").Value = "0"
Apply-DonorFormat "C28"

$ws.Range("M28").Value = 0
$ws.Range("N28").Value = -71.698113207547

# Row 29 - Shooting Inc.
Set-TextCell "C29"
$ws.Range("C29").Value = "0"
Apply-DonorFormat "C29"

$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -72.340425531914

Write-Host "edits applied"
